$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$ts = $sm.TextStyles
Write-Output ("TextStyles.Count=" + $ts.Count)
for ($i=1; $i -le $ts.Count; $i++) {
  $t = $ts.Item($i)
  $lvl1 = $t.Levels.Item(1)
  $f = $lvl1.Font
  Write-Output ("Style $i Level1 Font.Name=" + $f.Name)
}
